# This script applies a reshuffle of species-observation data among rows 34-40
# of the active worksheet:
#  - Rows 34/40, 35/39, 36/38 swap their Id/Taxonsorteringsordning/Rödlistade/
#    TaxonId/Artnamn/Vetenskapligt namn/Auktor (columns A,B,D,E,F,G,H) values.
#    Row 37 keeps its own values.
#  - The Ost/Nord (Q/R) coordinate values are rewritten as rounded integers
#    (taken from the row's new A/B/.../H source row) for every row 34-40.
#  - The Starttid/Sluttid (Z/AB) cells are cleared for rows 34-40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values to write into each row (34-40), already reflecting the swap
# described above plus the rounded Ost/Nord coordinates.
$rowData = @{
    34 = @{ A = 112038604; B = 89845;  D = "VU"; E = 1209; F = "Rynkskinn";        G = "Phlebia centrifuga";        H = "P.Karst.";                    Q = 615978; R = 6895550 }
    35 = @{ A = 112038601; B = 73634;  D = "LC"; E = 6426; F = "Kattfotslav";      G = "Felipes leucopellaeus";     H = "(Ach.) Frisch & G.Thor";      Q = 616013; R = 6895612 }
    36 = @{ A = 112038602; B = 86223;  D = "NT"; E = 4412; F = "Äggvaxskivling";   G = "Hygrophorus karstenii";     H = "Sacc. & Cub.";                Q = 616026; R = 6895554 }
    37 = @{                                                                                                                                            Q = 616034; R = 6895585 }
    38 = @{ A = 112038596; B = 90087;  D = "LC"; E = 3298; F = "Trådticka";        G = "Climacocystis borealis";    H = "(Fr.) Kotl. & Pouzar";        Q = 616076; R = 6895428 }
    39 = @{ A = 112038599; B = 89423;  D = "NT"; E = 5432; F = "Granticka";        G = "Porodaedalea chrysoloma";   H = "(Fr.) Fiasson & Niemelä";     Q = 616070; R = 6895500 }
    40 = @{ A = 112038603; B = 89369;  D = "LC"; E = 5447; F = "Vedticka";         G = "Fuscoporia viticola";       H = "(Schwein.) Murrill";          Q = 615968; R = 6895406 }
}

foreach ($r in 34..40) {
    $data = $rowData[$r]

    if ($data.ContainsKey("A")) { $ws.Cells.Item($r, 1).Value = $data.A }   # A: Id
    if ($data.ContainsKey("B")) { $ws.Cells.Item($r, 2).Value = $data.B }   # B: Taxonsorteringsordning
    if ($data.ContainsKey("D")) { $ws.Cells.Item($r, 4).Value = $data.D }   # D: Rödlistade
    if ($data.ContainsKey("E")) { $ws.Cells.Item($r, 5).Value = $data.E }   # E: TaxonId
    if ($data.ContainsKey("F")) { $ws.Cells.Item($r, 6).Value = $data.F }   # F: Artnamn
    if ($data.ContainsKey("G")) { $ws.Cells.Item($r, 7).Value = $data.G }   # G: Vetenskapligt namn
    if ($data.ContainsKey("H")) { $ws.Cells.Item($r, 8).Value = $data.H }   # H: Auktor

    $ws.Cells.Item($r, 17).Value = $data.Q   # Q: Ost
    $ws.Cells.Item($r, 18).Value = $data.R   # R: Nord

    $ws.Cells.Item($r, 26).ClearContents()   # Z: Starttid
    $ws.Cells.Item($r, 28).ClearContents()   # AB: Sluttid
}
